$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Fri Oct 25 16:00:02 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 16:00:20 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 16:00:38 EDT 2024"
$ws.Range("B5").Value = "Fri Oct 25 16:00:57 EDT 2024"
$ws.Range("B6").Value = "Fri Oct 25 16:01:15 EDT 2024"
$ws.Range("B7").Value = "Fri Oct 25 16:01:33 EDT 2024"

$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Fri Oct 25 16:01:52 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 16:02:10 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 16:02:28 EDT 2024"
$ws.Range("B5").Value = "Fri Oct 25 16:02:46 EDT 2024"
$ws.Range("B6").Value = "Fri Oct 25 16:03:04 EDT 2024"
$ws.Range("B7").Value = "Fri Oct 25 16:03:22 EDT 2024"
$ws.Range("B8").Value = "Fri Oct 25 16:03:42 EDT 2024"
$ws.Range("B9").Value = "Fri Oct 25 16:04:01 EDT 2024"
$ws.Range("B10").Value = "Fri Oct 25 16:04:20 EDT 2024"
$ws.Range("B11").Value = "Fri Oct 25 16:04:39 EDT 2024"
$ws.Range("B12").Value = "Fri Oct 25 16:04:56 EDT 2024"
$ws.Range("B13").Value = "Sat Oct 26 21:25:31 EDT 2024"
$ws.Range("B14").Value = "Sat Oct 26 21:25:52 EDT 2024"
$ws.Range("B15").Value = "Fri Oct 25 16:05:54 EDT 2024"
$ws.Range("B16").Value = "Fri Oct 25 16:06:12 EDT 2024"
$ws.Range("B17").Value = "Fri Oct 25 16:06:32 EDT 2024"
$ws.Range("B18").Value = "Fri Oct 25 16:06:51 EDT 2024"
$ws.Range("B19").Value = "Fri Oct 25 16:07:09 EDT 2024"

$ws = $wb.Worksheets.Item("Extension")
$ws.Range("B2").Value = "Fri Oct 25 16:07:28 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 16:07:46 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 16:08:04 EDT 2024"
$ws.Range("B5").Value = "Fri Oct 25 16:08:22 EDT 2024"
$ws.Range("B6").Value = "Fri Oct 25 16:08:40 EDT 2024"
$ws.Range("B7").Value = "Fri Oct 25 16:08:57 EDT 2024"

$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Fri Oct 25 16:09:15 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 16:09:34 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 16:09:52 EDT 2024"
$ws.Range("B5").Value = "Fri Oct 25 16:10:10 EDT 2024"
$ws.Range("B6").Value = "Fri Oct 25 16:10:29 EDT 2024"
$ws.Range("B7").Value = "Fri Oct 25 16:10:47 EDT 2024"
$ws.Range("B8").Value = "Fri Oct 25 16:11:06 EDT 2024"
$ws.Range("B9").Value = "Fri Oct 25 16:11:24 EDT 2024"
$ws.Range("B10").Value = "Fri Oct 25 16:11:43 EDT 2024"
$ws.Range("B11").Value = "Fri Oct 25 16:12:02 EDT 2024"
$ws.Range("B12").Value = "Fri Oct 25 16:12:20 EDT 2024"
$ws.Range("B13").Value = "Fri Oct 25 16:12:39 EDT 2024"
$ws.Range("B14").Value = "Fri Oct 25 16:12:57 EDT 2024"
$ws.Range("B15").Value = "Fri Oct 25 16:13:15 EDT 2024"
$ws.Range("B16").Value = "Fri Oct 25 16:13:34 EDT 2024"
$ws.Range("B17").Value = "Fri Oct 25 16:13:52 EDT 2024"
$ws.Range("B18").Value = "Fri Oct 25 16:14:13 EDT 2024"
$ws.Range("B19").Value = "Fri Oct 25 16:14:33 EDT 2024"
$ws.Range("B20").Value = "Fri Oct 25 16:14:51 EDT 2024"
$ws.Range("B21").Value = "Fri Oct 25 16:15:11 EDT 2024"
$ws.Range("B22").Value = "Fri Oct 25 16:15:31 EDT 2024"
$ws.Range("B23").Value = "Fri Oct 25 16:15:52 EDT 2024"
$ws.Range("B24").Value = "Fri Oct 25 16:16:10 EDT 2024"
$ws.Range("B25").Value = "Fri Oct 25 16:16:30 EDT 2024"
$ws.Range("B26").Value = "Fri Oct 25 16:16:50 EDT 2024"
$ws.Range("B27").Value = "Fri Oct 25 16:17:10 EDT 2024"
$ws.Range("B28").Value = "Fri Oct 25 16:17:29 EDT 2024"
$ws.Range("B29").Value = "Fri Oct 25 16:17:49 EDT 2024"
$ws.Range("B30").Value = "Fri Oct 25 16:18:09 EDT 2024"
$ws.Range("B31").Value = "Fri Oct 25 16:18:29 EDT 2024"
$ws.Range("B32").Value = "Fri Oct 25 16:18:47 EDT 2024"
$ws.Range("B33").Value = "Fri Oct 25 16:19:06 EDT 2024"
$ws.Range("B34").Value = "Fri Oct 25 16:19:27 EDT 2024"
$ws.Range("B35").Value = "Fri Oct 25 16:19:46 EDT 2024"
$ws.Range("B36").Value = "Fri Oct 25 16:20:05 EDT 2024"
$ws.Range("B37").Value = "Fri Oct 25 16:20:26 EDT 2024"
$ws.Range("B38").Value = "Fri Oct 25 16:20:46 EDT 2024"
$ws.Range("B39").Value = "Fri Oct 25 16:21:06 EDT 2024"
$ws.Range("B40").Value = "Fri Oct 25 16:21:24 EDT 2024"
$ws.Range("B41").Value = "Fri Oct 25 16:21:44 EDT 2024"
$ws.Range("B42").Value = "Fri Oct 25 16:22:04 EDT 2024"
$ws.Range("B43").Value = "Fri Oct 25 16:22:25 EDT 2024"
$ws.Range("B44").Value = "Fri Oct 25 16:22:43 EDT 2024"
$ws.Range("B45").Value = "Fri Oct 25 16:23:01 EDT 2024"
$ws.Range("B46").Value = "Fri Oct 25 16:23:22 EDT 2024"
$ws.Range("B47").Value = "Fri Oct 25 16:23:43 EDT 2024"
$ws.Range("B48").Value = "Fri Oct 25 16:24:03 EDT 2024"
$ws.Range("B49").Value = "Fri Oct 25 16:24:23 EDT 2024"
$ws.Range("B50").Value = "Fri Oct 25 16:24:43 EDT 2024"
$ws.Range("B51").Value = "Fri Oct 25 16:25:03 EDT 2024"
$ws.Range("B52").Value = "Fri Oct 25 16:25:21 EDT 2024"

$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Fri Oct 25 16:26:22 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 16:26:40 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 16:26:58 EDT 2024"
$ws.Range("B5").Value = "Fri Oct 25 16:27:17 EDT 2024"
$ws.Range("B6").Value = "Fri Oct 25 16:27:35 EDT 2024"
$ws.Range("B7").Value = "Fri Oct 25 16:27:53 EDT 2024"
$ws.Range("B8").Value = "Fri Oct 25 16:28:13 EDT 2024"
$ws.Range("B9").Value = "Fri Oct 25 16:28:31 EDT 2024"

$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Fri Oct 25 16:28:51 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 16:29:16 EDT 2024"
$ws.Range("B4").Value = "Fri Oct 25 16:29:41 EDT 2024"
$ws.Range("B5").Value = "Fri Oct 25 16:30:06 EDT 2024"
$ws.Range("B6").Value = "Fri Oct 25 16:30:32 EDT 2024"

$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Fri Oct 25 16:25:42 EDT 2024"
$ws.Range("B3").Value = "Fri Oct 25 16:26:01 EDT 2024"
